# The sheet's data rows (A2:A66) were originally stamped with the bold /
# bordered / centered "header" style (style index 1) along with the real
# header row (A1:B1). This clears that formatting back to the workbook's
# default ("Normal") style on the data rows only, leaving the header row
# (A1:B1) with its original header style intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A66").Style = "Normal"

# Re-affirm the header labels (no visible change - they already hold these
# values) so the header row keeps its original "Announced Date" /
# "Money Raised Currency (in USD)" text and style.
$ws.Range("A1").Value = "Announced Date"
$ws.Range("B1").Value = "Money Raised Currency (in USD)"
